# 07_abstract_classes_interfaces.pptx
# Commit: "corrected wrong lasson enum"
#
# The title slide's subtitle still announced this deck as lesson "06" while
# the rest of the course (filename, content) is lesson 07. Fix the stray
# lesson-number typo in the subtitle placeholder on the first slide.

$p = $ppt.ActivePresentation

$oldText = "06 - Abstrakte Klassen und Interfaces"
$newText = "07 - Abstrakte Klassen und Interfaces"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}
